$d = $word.ActiveDocument

function Get-ContainingParagraph($doc, $startPos, $endPos) {
    $found = $null
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Start -le $startPos -and $endPos -le $p.Range.End) {
            $found = $p
        }
    }
    return $found
}

# ------------------------------------------------------------------
# 1) Replace the "m:if self.name <> 'anydsl'" field (fldChar/instrText)
#    with plain text runs "{m:if self.name <> 'anydsl'}" split across
#    7 runs, matching the TokenIteratorFieldRewriterSplit output.
# ------------------------------------------------------------------
$ifField = $null
foreach ($f in $d.Fields) {
    if ($f.Code.Text.Trim().StartsWith("m:if")) {
        $ifField = $f
        break
    }
}
$ifParagraph = Get-ContainingParagraph $d $ifField.Code.Start $ifField.Code.End
$ifRange = $d.Range($ifParagraph.Range.Start, $ifParagraph.Range.End - 1)

$ifXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="002A1F2A" w:rsidRDefault="002A1F2A" w:rsidP="002A1F2A"><w:pPr><w:tabs><w:tab w:val="left" w:pos="3119"/></w:tabs></w:pPr><w:r><w:t xml:space="preserve">{m:if </w:t></w:r><w:r><w:t xml:space="preserve">self.name </w:t></w:r><w:r><w:t>&lt;&gt;</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>'</w:t></w:r><w:r><w:t>anydsl</w:t></w:r><w:r><w:t>'}</w:t></w:r></w:p></w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@
[void]$ifRange.InsertXML($ifXml)

# ------------------------------------------------------------------
# 2) Replace the "m:else" field (fldChar/instrText) with a single
#    plain text run "{m:else}".
# ------------------------------------------------------------------
$elseField = $null
foreach ($f in $d.Fields) {
    if ($f.Code.Text.Trim().StartsWith("m:else")) {
        $elseField = $f
        break
    }
}
$elseParagraph = Get-ContainingParagraph $d $elseField.Code.Start $elseField.Code.End
$elseRange = $d.Range($elseParagraph.Range.Start, $elseParagraph.Range.End - 1)

$elseXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00A07687" w:rsidRDefault="00A07687" w:rsidP="00A07687"><w:pPr><w:tabs><w:tab w:val="left" w:pos="3119"/></w:tabs></w:pPr><w:r><w:t>{m:else}</w:t></w:r></w:p></w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@
[void]$elseRange.InsertXML($elseXml)

# ------------------------------------------------------------------
# 3) Remove the trailing "    " (4 spaces) run that follows the
#    "Unexpected tag EOF ..." error message run in the last paragraph.
#    (Delete only the trailing whitespace characters/run, leaving the
#    preceding error-message run untouched and separate.)
# ------------------------------------------------------------------
$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)
$fullText = $lastParagraph.Range.Text
$withoutMark = $fullText.Substring(0, $fullText.Length - 1)
$trimmed = $withoutMark.TrimEnd(" ")
if ($trimmed.Length -lt $withoutMark.Length) {
    $deleteStart = $lastParagraph.Range.Start + $trimmed.Length
    $deleteEnd = $lastParagraph.Range.Start + $withoutMark.Length
    $trailingRange = $d.Range($deleteStart, $deleteEnd)
    [void]$trailingRange.Delete()
}
